# Update column F (dSF) values after repulling data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = -3
    4  = 2
    5  = 1
    6  = 2
    7  = -3
    8  = -2
    9  = 2
    10 = -2
    12 = 2
    13 = 4
    14 = -1
    16 = 1
    17 = 2
    18 = -2
    20 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
